# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) across the leve-profit sheets with refreshed
# Universalis price snapshots. Values only; no formulas involved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 51.42857
$ws.Range("I8").Value = 51.42857
$ws.Range("K8").Value = 154.28571
$ws.Range("M8").Value = -15.28570999999999
$ws.Range("H43").Value = 2428.111
$ws.Range("J43").Value = 2994.5
$ws.Range("L43").Value = 2994.5
$ws.Range("N43").Value = -3132.5
$ws.Range("H64").Value = 7749.231
$ws.Range("I64").Value = 3656.6667
$ws.Range("K64").Value = 3656.6667
$ws.Range("M64").Value = -3408.6667
$ws.Range("H67").Value = 7749.231
$ws.Range("I67").Value = 3656.6667
$ws.Range("K67").Value = 3656.6667
$ws.Range("M67").Value = -2798.6667
$ws.Range("H70").Value = 4864.6665
$ws.Range("I70").Value = 4200
$ws.Range("J70").Value = 4997.6
$ws.Range("K70").Value = 12600
$ws.Range("L70").Value = 14992.8
$ws.Range("M70").Value = -12330
$ws.Range("N70").Value = -15532.8
$ws.Range("H73").Value = 4864.6665
$ws.Range("I73").Value = 4200
$ws.Range("J73").Value = 4997.6
$ws.Range("K73").Value = 12600
$ws.Range("L73").Value = 14992.8
$ws.Range("M73").Value = -11664
$ws.Range("N73").Value = -16864.8
$ws.Range("H74").Value = 6558.1055
$ws.Range("I74").Value = 3686
$ws.Range("K74").Value = 3686
$ws.Range("M74").Value = -2750
$ws.Range("H77").Value = 6558.1055
$ws.Range("I77").Value = 3686
$ws.Range("K77").Value = 18430
$ws.Range("M77").Value = -13750
$ws.Range("H113").Value = 1837
$ws.Range("I113").Value = 1837
$ws.Range("K113").Value = 1837
$ws.Range("M113").Value = 1417
$ws.Range("H121").Value = 7898.7144
$ws.Range("I121").Value = 7898.7144
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 23696.1432
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -21949.1432
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 1698.25
$ws.Range("I132").Value = 1625.3636
$ws.Range("K132").Value = 4876.0908
$ws.Range("M132").Value = -2346.0908
$ws.Range("H137").Value = 1757.4117
$ws.Range("J137").Value = 2183
$ws.Range("L137").Value = 6549
$ws.Range("N137").Value = -11649
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1804.7216
$ws.Range("I32").Value = 1868.1184
$ws.Range("K32").Value = 1868.1184
$ws.Range("M32").Value = -1581.1184
$ws.Range("H55").Value = 32200
$ws.Range("J55").Value = 38000
$ws.Range("L55").Value = 38000
$ws.Range("N55").Value = -38630
$ws.Range("H74").Value = 2878.9697
$ws.Range("I74").Value = 1919.3572
$ws.Range("J74").Value = 8252.799999999999
$ws.Range("K74").Value = 1919.3572
$ws.Range("L74").Value = 8252.799999999999
$ws.Range("M74").Value = -1045.3572
$ws.Range("N74").Value = -10000.8
$ws.Range("H77").Value = 2878.9697
$ws.Range("I77").Value = 1919.3572
$ws.Range("J77").Value = 8252.799999999999
$ws.Range("K77").Value = 9596.786
$ws.Range("L77").Value = 41264
$ws.Range("M77").Value = -5228.786
$ws.Range("N77").Value = -50000
$ws.Range("H80").Value = 41300
$ws.Range("J80").Value = 41300
$ws.Range("L80").Value = 41300
$ws.Range("N80").Value = -43296
$ws.Range("H83").Value = 41300
$ws.Range("J83").Value = 41300
$ws.Range("L83").Value = 123900
$ws.Range("N83").Value = -133884
$ws.Range("H132").Value = 4358.8096
$ws.Range("I132").Value = 3846.5789
$ws.Range("K132").Value = 11539.7367
$ws.Range("M132").Value = -9009.736699999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1593.3684
$ws.Range("I105").Value = 1602
$ws.Range("K105").Value = 1602
$ws.Range("M105").Value = 145
$ws.Range("H134").Value = 2615.175
$ws.Range("J134").Value = 1616.1666
$ws.Range("L134").Value = 4848.4998
$ws.Range("N134").Value = -9918.4998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4380.5884
$ws.Range("I31").Value = 3452.5557
$ws.Range("J31").Value = 5424.625
$ws.Range("K31").Value = 3452.5557
$ws.Range("L31").Value = 5424.625
$ws.Range("M31").Value = -3157.5557
$ws.Range("N31").Value = -6014.625
$ws.Range("H34").Value = 4380.5884
$ws.Range("I34").Value = 3452.5557
$ws.Range("J34").Value = 5424.625
$ws.Range("K34").Value = 3452.5557
$ws.Range("L34").Value = 5424.625
$ws.Range("M34").Value = -3250.5557
$ws.Range("N34").Value = -5828.625
$ws.Range("H41").Value = 13753.333
$ws.Range("J41").Value = 19513.334
$ws.Range("L41").Value = 19513.334
$ws.Range("N41").Value = -20369.334
$ws.Range("H58").Value = 7556.737
$ws.Range("I58").Value = 3154
$ws.Range("K58").Value = 3154
$ws.Range("M58").Value = -2951
$ws.Range("H122").Value = 4389.6665
$ws.Range("I122").Value = 4559.0713
$ws.Range("K122").Value = 13677.2139
$ws.Range("M122").Value = -11227.2139
$ws.Range("H131").Value = 37999
$ws.Range("J131").Value = 39499.5
$ws.Range("L131").Value = 39499.5
$ws.Range("N131").Value = -49579.5
$ws.Range("H132").Value = 2868.8333
$ws.Range("I132").Value = 2369.5557
$ws.Range("J132").Value = 4366.6665
$ws.Range("K132").Value = 7108.6671
$ws.Range("L132").Value = 13099.9995
$ws.Range("M132").Value = -4578.6671
$ws.Range("N132").Value = -18159.9995
$ws.Range("H136").Value = 7556.737
$ws.Range("I136").Value = 3154
$ws.Range("K136").Value = 9462
$ws.Range("M136").Value = -6912
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4708869
$ws.Range("I4").Value = 1465651.8
$ws.Range("K4").Value = 4396955.4
$ws.Range("M4").Value = -4396843.4
$ws.Range("H8").Value = 1232
$ws.Range("I8").Value = 1232
$ws.Range("K8").Value = 3696
$ws.Range("M8").Value = -3557
$ws.Range("H38").Value = 468.73685
$ws.Range("J38").Value = 1102.8572
$ws.Range("L38").Value = 3308.5716
$ws.Range("N38").Value = -4002.5716
$ws.Range("H99").Value = 1750
$ws.Range("I99").Value = 1750
$ws.Range("K99").Value = 5250
$ws.Range("M99").Value = -3004
$ws.Range("H105").Value = 14398
$ws.Range("J105").Value = 14398
$ws.Range("L105").Value = 43194
$ws.Range("N105").Value = -48436
$ws.Range("H106").Value = 4160.778
$ws.Range("I106").Value = 3222
$ws.Range("J106").Value = 4278.125
$ws.Range("K106").Value = 9666
$ws.Range("L106").Value = 12834.375
$ws.Range("M106").Value = -8720
$ws.Range("N106").Value = -14726.375
$ws.Range("H113").Value = 1566.6666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1566.6666
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4699.9998
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9039.9998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 50000
$ws.Range("I26").Value = 50000
$ws.Range("K26").Value = 50000
$ws.Range("M26").Value = -49720
$ws.Range("H50").Value = 50000
$ws.Range("I50").Value = 50000
$ws.Range("K50").Value = 50000
$ws.Range("M50").Value = -49502
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31262
$ws.Range("H126").Value = 2712.9443
$ws.Range("J126").Value = 3874.875
$ws.Range("L126").Value = 11624.625
$ws.Range("N126").Value = -16564.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6466.778
$ws.Range("I7").Value = 2481.4285
$ws.Range("K7").Value = 2481.4285
$ws.Range("M7").Value = -2369.4285
$ws.Range("H22").Value = 1899.8889
$ws.Range("J22").Value = 2333.1667
$ws.Range("L22").Value = 2333.1667
$ws.Range("N22").Value = -2923.1667
$ws.Range("H27").Value = 1899.8889
$ws.Range("J27").Value = 2333.1667
$ws.Range("L27").Value = 2333.1667
$ws.Range("N27").Value = -2547.1667
$ws.Range("H55").Value = 879.375
$ws.Range("I55").Value = 1212.5714
$ws.Range("J55").Value = 620.2222
$ws.Range("K55").Value = 1212.5714
$ws.Range("L55").Value = 620.2222
$ws.Range("M55").Value = -1039.5714
$ws.Range("N55").Value = -966.2222
$ws.Range("H68").Value = 13040.63
$ws.Range("I68").Value = 13436.546
$ws.Range("J68").Value = 11298.6
$ws.Range("K68").Value = 13436.546
$ws.Range("L68").Value = 11298.6
$ws.Range("M68").Value = -12687.546
$ws.Range("N68").Value = -12796.6
$ws.Range("H71").Value = 13040.63
$ws.Range("I71").Value = 13436.546
$ws.Range("J71").Value = 11298.6
$ws.Range("K71").Value = 67182.73
$ws.Range("L71").Value = 56493
$ws.Range("M71").Value = -63438.73
$ws.Range("N71").Value = -63981
$ws.Range("H126").Value = 6466.778
$ws.Range("I126").Value = 2481.4285
$ws.Range("K126").Value = 7444.2855
$ws.Range("M126").Value = -4974.2855
$ws.Range("H132").Value = 6130.625
$ws.Range("I132").Value = 6130.625
$ws.Range("K132").Value = 18391.875
$ws.Range("M132").Value = -15861.875
$ws.Range("H136").Value = 8181.5454
$ws.Range("I136").Value = 7433.4287
$ws.Range("K136").Value = 22300.2861
$ws.Range("M136").Value = -19750.2861
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3159
$ws.Range("I132").Value = 2300.2222
$ws.Range("K132").Value = 6900.6666
$ws.Range("M132").Value = -4370.6666
$ws.Range("H136").Value = 5863.636
$ws.Range("I136").Value = 4950
$ws.Range("K136").Value = 14850
$ws.Range("M136").Value = -12300
